$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capitalize first letter of the "customer/beta testing" row header (F3)
$ws.Range("F3").Value = "Testovanie zákazníkom a Beta testy"

# Update the view: scroll so column B is the left-most visible column, and select F3
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F3").Select()
